$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.440.89'
$ws.Range("E2").Value = '  +3.55%  '
$ws.Range("D3").Value = '1.589.52'
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("E4").Value = '  +0.95%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.492'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.66%  '
$ws.Range("E7").Value = '  +0.90%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.45'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +8.45%  '
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("E10").Value = '  +1.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0887'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.97%  '
$ws.Range("D12").Value = '1.816.72'
$ws.Range("E12").Value = '  +1.56%  '
$ws.Range("D13").Value = '1.588.85'
$ws.Range("E13").Value = '  +1.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.531'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.75'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.24%  '
$ws.Range("D16").Value = '28.459.02'
$ws.Range("E16").Value = '  +3.54%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.09'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '230.14'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.53%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.0₃0707'
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("E21").Value = '  +1.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.34'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.68%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.57'
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = '  +0.63%  '
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("E29").Value = '  +0.91%  '
$ws.Range("E30").Value = '  -1.00%  '
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("E32").Value = '  +0.64%  '
$ws.Range("E33").Value = '  +1.02%  '
$ws.Range("D34").Value = '1.400.38'
$ws.Range("E34").Value = '  -3.10%  '
$ws.Range("E35").Value = '  -0.80%  '
$ws.Range("E36").Value = '  -9.07%  '
$ws.Range("E37").Value = '  +1.03%  '
$ws.Range("E38").Value = '  -0.19%  '
$ws.Range("E39").Value = '  +7.72%  '
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("E42").Value = '  +0.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.59'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.47%  '
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.979'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '63.57'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.22%  '
$ws.Range("D47").Value = '1.726.36'
$ws.Range("E47").Value = '  +1.38%  '
$ws.Range("B48").Value = 'mCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.86%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '87.36'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.86%  '
$ws.Range("D50").Value = '0.0₆0105'
$ws.Range("E50").Value = '  +8.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0520'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.03%  '
